$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date and the concept Count ---
$meta = $wb.Worksheets.Item("Metadata")

# Date is plain (non-numeric-looking) text, safe to set directly.
$meta.Range("B8").Value = "2022-11-17T12:54:43-06:00"

# Count must stay a shared-string "2" (not a numeric cell), matching how the
# sheet already stores "1" as text. Seed the cell with a same-style text cell,
# write the new digit through a text formula (so Excel can't re-infer it as a
# number), then flatten the formula back down to a plain value in place.
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A2").Copy()
$meta.Range("B23").PasteSpecial(-4104)
$meta.Range("B23").Formula = "=""2"""
$meta.Range("B23").Copy()
$meta.Range("B23").PasteSpecial(-4163)

# --- Concepts sheet: add a second, explicit Apple ("Sync My Health Data") row ---

# Duplicate row 2's formatting + values down into row 3 (keeps style index and
# the Level cell "1" as the same shared text value as row 2's A2).
$concepts.Range("A2:D2").Copy()
$concepts.Range("A3:D3").PasteSpecial(-4104)
$concepts.Range("A2:D2").Copy()
$concepts.Range("A3:D3").PasteSpecial(-4122)

# New concept: second top-level ("Level" 1) entry, representing the Apple app.
$concepts.Range("B3").Value = "49d985b4-6f9d-47e3-bbc9-a56c840ba4a6"
$concepts.Range("C3").Value = "Sync My Health Data"
